$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect to allow edits, then restore protection.
$ws.Unprotect("D382")

# Update the "as of" date in the confidential disclaimer text (A59).
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-03 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for holdings rows 2-56.
$ws.Range("D2").Value = 0.01101565450621742
$ws.Range("E2").Value = 0.03439716312056751
$ws.Range("D3").Value = 0.009898017583975668
$ws.Range("E3").Value = 0.01968743657398009
$ws.Range("D4").Value = 0.01074110019937553
$ws.Range("E4").Value = 0.01795511221945145
$ws.Range("D5").Value = 0.01140014214323788
$ws.Range("E5").Value = 0.01790591805766306
$ws.Range("D6").Value = 0.01084288618630228
$ws.Range("E6").Value = 0.02519762845849804
$ws.Range("D7").Value = 0.01212715301915987
$ws.Range("E7").Value = 0.01243339253996467
$ws.Range("D8").Value = 0.01137570011348244
$ws.Range("E8").Value = -0.007240547063555924
$ws.Range("D9").Value = 0.01116877991637477
$ws.Range("E9").Value = 0.02110480454073071
$ws.Range("D10").Value = 0.01046654589253364
$ws.Range("E10").Value = 0.02725527831094054
$ws.Range("D11").Value = 0.01097212760391322
$ws.Range("E11").Value = 0.02136100091547122
$ws.Range("D12").Value = 0.4506373677786636
$ws.Range("E12").Value = 0
$ws.Range("D13").Value = 0.011983514241556
$ws.Range("E13").Value = 0.007283090563647843
$ws.Range("D14").Value = 0.01099377944762352
$ws.Range("E14").Value = 0.009096077316657292
$ws.Range("D15").Value = 0.01052368890273813
$ws.Range("E15").Value = -0.001399906672888473
$ws.Range("D16").Value = 0.01017078617174867
$ws.Range("E16").Value = 0.003423680456490885
$ws.Range("D17").Value = 0.01012223693456321
$ws.Range("E17").Value = -0.02031534263189794
$ws.Range("D18").Value = 0.008137298582049771
$ws.Range("E18").Value = 0.06119873817034693
$ws.Range("D19").Value = 0.008693215249683703
$ws.Range("E19").Value = 0.04067222143764981
$ws.Range("D20").Value = 0.01268976613330214
$ws.Range("E20").Value = 0.009498680738786236
$ws.Range("D21").Value = 0.01178440656537472
$ws.Range("E21").Value = 0.0144902829867033
$ws.Range("D22").Value = 0.01158195066593928
$ws.Range("E22").Value = 0.04116638078902235
$ws.Range("D23").Value = 0.01181755397559101
$ws.Range("E23").Value = 0.005949851253718608
$ws.Range("D24").Value = 0.01187391573370286
$ws.Range("E24").Value = 0.01409906946141559
$ws.Range("D25").Value = 0.01245003334830361
$ws.Range("E25").Value = 0.01391279403327594
$ws.Range("D26").Value = 0.01110136902152416
$ws.Range("E26").Value = 0.03794185064543365
$ws.Range("D27").Value = 0.010303822242381
$ws.Range("E27").Value = 0.03338315894369703
$ws.Range("D28").Value = 0.01218719782285131
$ws.Range("E28").Value = 0.02085222121486852
$ws.Range("D29").Value = 0.0102623042740293
$ws.Range("E29").Value = 0.04469820554649284
$ws.Range("D30").Value = 0.006856938009655381
$ws.Range("E30").Value = 0.02031316123571725
$ws.Range("D31").Value = 0.00529298292763284
$ws.Range("E31").Value = -0.01992619926199268
$ws.Range("D32").Value = 0.009088752023442918
$ws.Range("E32").Value = 0.02404371584699438
$ws.Range("D33").Value = 0.01077145742354667
$ws.Range("E33").Value = 0.03545673076923084
$ws.Range("D34").Value = 0.01048094325252657
$ws.Range("E34").Value = 0.02792064658339455
$ws.Range("D35").Value = 0.009456498544192527
$ws.Range("E35").Value = 0.01746724890829698
$ws.Range("D36").Value = 0.01022033987591038
$ws.Range("E36").Value = -0.009478672985781977
$ws.Range("D37").Value = 0.009904044385833173
$ws.Range("E37").Value = 0.005172413793103292
$ws.Range("D38").Value = 0.01141130288741845
$ws.Range("E38").Value = 0.0250134480903712
$ws.Range("D39").Value = 0.01333575000647323
$ws.Range("E39").Value = 0.03093197643277978
$ws.Range("D40").Value = 0.0112049407275198
$ws.Range("E40").Value = 0.03557910673732034
$ws.Range("D41").Value = 0.01173574572074746
$ws.Range("E41").Value = 0.02967133292757151
$ws.Range("D42").Value = 0.01131308833862948
$ws.Range("E42").Value = -0.002249297094658043
$ws.Range("D43").Value = 0.01118094512753159
$ws.Range("E43").Value = 0.0255138199858258
$ws.Range("D44").Value = 0.01089935955185594
$ws.Range("E44").Value = -0.009185115402732036
$ws.Range("D45").Value = 0.01180248697094724
$ws.Range("E45").Value = 0.0271867612293144
$ws.Range("D46").Value = 0.01080851109422614
$ws.Range("E46").Value = 0.02044525215810999
$ws.Range("D47").Value = 0.01023283990939261
$ws.Range("E47").Value = 0.01282638570774175
$ws.Range("D48").Value = 0.01047324233904198
$ws.Range("E48").Value = 0.03708439897698201
$ws.Range("D49").Value = 0.009780160125428905
$ws.Range("E49").Value = 0.01889763779527565
$ws.Range("D50").Value = 0.009642882972007959
$ws.Range("E50").Value = 0.03916666666666679
$ws.Range("D51").Value = 0.009403485009334845
$ws.Range("E51").Value = -0.001459854014598694
$ws.Range("D52").Value = 0.0100290447206555
$ws.Range("E52").Value = 0.00689962163365232
$ws.Range("D53").Value = 0.009099577945298067
$ws.Range("E53").Value = 0.01084236864053367
$ws.Range("D54").Value = 0.00414342627703467
$ws.Range("E54").Value = 0.03313131313131312
$ws.Range("D55").Value = 0.004108939577516724
$ws.Range("E55").Value = 0
$ws.Range("E56").Value = 0.01014444681548232

$ws.Protect("D382")
